# Updated for Insurance Tables
# Adds a new "InsuranceTables" worksheet (positioned before "RegionWeights"),
# and refreshes the "RegionWeights" sheet view state (it moves from being the
# tab-selected sheet to a background sheet, scrolled further down).

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "InsuranceTables" sheet.
#    Worksheets.Add() on this host inserts right before the sheet that
#    used to be last ("RegionWeights"), which is exactly where the diff
#    wants it (sheetId 12, rId10 reused, RegionWeights bumped to rId11).
# ------------------------------------------------------------------
$ins = $wb.Worksheets.Add()
$ins.Name = "InsuranceTables"

# ------------------------------------------------------------------
# 2. Header row.
# ------------------------------------------------------------------
$ins.Range("A1").Value = "Year"
$ins.Range("B1").Value = "Table"
$ins.Range("C1").Value = "StartCode"
$ins.Range("D1").Value = "EndCode"
$ins.Range("E1").Value = "HHID"
$ins.Range("F1").Value = "Code"
$ins.Range("G1").Value = "InsuredCount"
$ins.Range("H1").Value = "InsuranceCosts"
$ins.Range("I1").Value = "GovPaid"
$ins.Range("J1").Value = "SS1"
$ins.Range("K1").Value = "SS2"

# ------------------------------------------------------------------
# 3. Column A year codes, rows 2-28 (63..89) - only column A populated.
# ------------------------------------------------------------------
$year = 63
for ($r = 2; $r -le 28; $r++) {
    $ins.Cells.Item($r, 1).Value = $year
    $year = $year + 1
}

# ------------------------------------------------------------------
# 4. Fully populated rows 29-33 (years 90..94), identical row pattern.
# ------------------------------------------------------------------
$year = 90
for ($r = 29; $r -le 33; $r++) {
    $ins.Cells.Item($r, 1).Value = $year
    $ins.Cells.Item($r, 2).Value = "P3S13"
    $ins.Cells.Item($r, 3).Value = 125311
    $ins.Cells.Item($r, 4).Value = 125317
    $ins.Cells.Item($r, 5).Value = "Address"
    $ins.Cells.Item($r, 6).Value = "DYCOL01"
    $ins.Cells.Item($r, 7).Value = "DYCOL02"
    $ins.Cells.Item($r, 8).Value = "DYCOL05"
    $ins.Cells.Item($r, 9).Value = 125312
    $ins.Cells.Item($r, 10).Value = 125313
    $ins.Cells.Item($r, 11).Value = 125314
    $year = $year + 1
}

# ------------------------------------------------------------------
# 5. Column H width (13.875 in the target file).
# ------------------------------------------------------------------
$ins.Columns.Item(8).ColumnWidth = 13.1

# ------------------------------------------------------------------
# 6. View: freeze header row + first column (split at B2 -> xSplit=1,
#    ySplit=1), then move the bottom-right selection to H33, and make
#    this the tab-selected sheet (it becomes the active tab).
# ------------------------------------------------------------------
$ins.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ins.Range("H33").Select()
$ins.Activate()

# ------------------------------------------------------------------
# 7. RegionWeights view refresh: no longer the tab-selected sheet, pane
#    scrolled further down, and the bottom-right selection moved.
#    (Re-freezing at B2 keeps xSplit=1 / ySplit=1 as before.)
# ------------------------------------------------------------------
$rw = $wb.Worksheets.Item("RegionWeights")
$rw.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$rw.Range("G113").Select()

# Leave InsuranceTables as the active/selected sheet (matches tabSelected
# moving onto the new sheet).
$ins.Activate()
